# Apply calibrated value updates to rows 96-101, 103-104, 111-112
# (columns J:AS, i.e. the full time-series of each row, all set to the
# same new constant value per row) on the single worksheet "strategy_id-0".

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("strategy_id-0")

$rowValues = @{
    96  = 23403862.31
    97  = 17500919.01
    98  = 4452377.762
    99  = 296129.2898
    100 = 3150919.414
    101 = 28032991.61
    103 = 4939215.276
    104 = 4635387.512
    111 = 281075.8373
    112 = 777959.6513
}

# Columns J (10) through AS (45) hold the 36 year-series values (years 0-35)
# for each row; every one of those cells gets set to the same new constant.
$numCols = 36

foreach ($row in $rowValues.Keys) {
    $value = $rowValues[$row]
    $rng = $ws.Range("J$row").Resize(1, $numCols)
    $arr = New-Object 'object[,]' 1,$numCols
    for ($i = 0; $i -lt $numCols; $i++) {
        $arr[0,$i] = $value
    }
    $rng.Value = $arr
}

$wb.Save()
